# RPA datasets push 2024-07-16
# Refresh the IPO underwriting table on Sheet1 with the latest RPA pull:
#  - one new listing ("엑셀세라퓨틱스", 대신 underwriter) is inserted as row 10
#  - every other row shifts/refreshes to match the new source order & values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert the new row so the table grows from 28 to 29 data rows (A1:L28 -> A1:L29).
$ws.Rows.Item(10).Insert()

# Row 2: 에스오에스랩
$ws.Range("A2").Value = "BNK"
$ws.Range("B2").Formula = '="2024-06-14"'
$ws.Range("C2").Value = "에스오에스랩"
$ws.Range("D2").Value = "한국"
$ws.Range("E2").Value = "한국, BNK"
$ws.Range("F2").Formula = '="2024-06-19"'
$ws.Range("G2").Formula = '="2024-06-25"'
$ws.Range("H2").Value = 1150
$ws.Range("I2").Value = 2000000
$ws.Range("J2").Value = 11500
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 5

# Row 3: 디비금융스팩12호
$ws.Range("A3").Value = "DB"
$ws.Range("B3").Formula = '="2024-06-05"'
$ws.Range("C3").Value = "디비금융스팩12호"
$ws.Range("D3").Value = "DB"
$ws.Range("E3").Value = "DB"
$ws.Range("F3").Formula = '="2024-06-11"'
$ws.Range("G3").Formula = '="2024-06-18"'
$ws.Range("H3").Value = 10000
$ws.Range("I3").Value = 5000000
$ws.Range("J3").Value = 2000
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 100

# Row 4: 한중엔시에스
$ws.Range("A4").Value = "IBK"
$ws.Range("B4").Formula = '="2024-06-10"'
$ws.Range("C4").Value = "한중엔시에스"
$ws.Range("D4").Value = "IBK"
$ws.Range("E4").Value = "IBK"
$ws.Range("F4").Formula = '="2024-06-13"'
$ws.Range("G4").Formula = '="2024-06-24"'
$ws.Range("H4").Value = 48000
$ws.Range("I4").Value = 1600000
$ws.Range("J4").Value = 30000
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 100

# Row 5: KB제29호스팩
$ws.Range("A5").Value = "KB"
$ws.Range("B5").Formula = '="2024-06-11"'
$ws.Range("C5").Value = "KB제29호스팩"
$ws.Range("D5").Value = "KB"
$ws.Range("E5").Value = "KB"
$ws.Range("F5").Formula = '="2024-06-14"'
$ws.Range("G5").Formula = '="2024-06-21"'
$ws.Range("H5").Value = 12000
$ws.Range("I5").Value = 6000000
$ws.Range("J5").Value = 2000
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 100

# Row 6: 이베스트스팩6호
$ws.Range("A6").Value = "LS"
$ws.Range("B6").Formula = '="2024-07-02"'
$ws.Range("C6").Value = "이베스트스팩6호"
$ws.Range("D6").Value = "엘에스"
$ws.Range("E6").Value = "엘에스"
$ws.Range("F6").Formula = '="2024-07-05"'
$ws.Range("G6").Formula = '="2024-07-12"'
$ws.Range("H6").Value = 8000
$ws.Range("I6").Value = 4000000
$ws.Range("J6").Value = 2000
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 100

# Row 7: 시프트업
$ws.Range("A7").Value = "NH"
$ws.Range("B7").Formula = '="2024-07-02"'
$ws.Range("C7").Value = "시프트업"
$ws.Range("D7").Value = "한국, 제이피모간회사 서울지점, NH"
$ws.Range("E7").Value = "한국, 제이피모간회사 서울지점, NH, 신한"
$ws.Range("F7").Formula = '="2024-07-05"'
$ws.Range("G7").Formula = '="2024-07-11"'
$ws.Range("H7").Value = 130500
$ws.Range("I7").Value = 7250000
$ws.Range("J7").Value = 60000
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 30

# Row 8: 에이치브이엠
$ws.Range("A8").Value = "NH"
$ws.Range("B8").Formula = '="2024-06-19"'
$ws.Range("C8").Value = "에이치브이엠"
$ws.Range("D8").Value = "NH"
$ws.Range("E8").Value = "NH"
$ws.Range("F8").Formula = '="2024-06-24"'
$ws.Range("G8").Formula = '="2024-06-28"'
$ws.Range("H8").Value = 43200
$ws.Range("I8").Value = 2400000
$ws.Range("J8").Value = 18000
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 100

# Row 9: 라메디텍
$ws.Range("A9").Value = "대신"
$ws.Range("B9").Formula = '="2024-06-05"'
$ws.Range("C9").Value = "라메디텍"
$ws.Range("D9").Value = "대신"
$ws.Range("E9").Value = "대신"
$ws.Range("F9").Formula = '="2024-06-11"'
$ws.Range("G9").Formula = '="2024-06-17"'
$ws.Range("H9").Value = 20768
$ws.Range("I9").Value = 1298000
$ws.Range("J9").Value = 16000
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 100

# Row 10: 엑셀세라퓨틱스
$ws.Range("A10").Value = "대신"
$ws.Range("B10").Formula = '="2024-07-03"'
$ws.Range("C10").Value = "엑셀세라퓨틱스"
$ws.Range("D10").Value = "대신"
$ws.Range("E10").Value = "대신"
$ws.Range("F10").Formula = '="2024-07-08"'
$ws.Range("G10").Formula = '="2024-07-15"'
$ws.Range("H10").Value = 16180
$ws.Range("I10").Value = 1618000
$ws.Range("J10").Value = 10000
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 100

# Row 11: 이노스페이스
$ws.Range("A11").Value = "미래"
$ws.Range("B11").Formula = '="2024-06-20"'
$ws.Range("C11").Value = "이노스페이스"
$ws.Range("D11").Value = "미래"
$ws.Range("E11").Value = "미래, 신한"
$ws.Range("F11").Formula = '="2024-06-25"'
$ws.Range("G11").Formula = '="2024-07-02"'
$ws.Range("H11").Value = 54133.66
$ws.Range("I11").Value = 1330000
$ws.Range("J11").Value = 43300
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 94

# Row 12: 미래에셋비전스팩4호
$ws.Range("A12").Value = "미래"
$ws.Range("B12").Formula = '="2024-05-20"'
$ws.Range("C12").Value = "미래에셋비전스팩4호"
$ws.Range("D12").Value = "미래"
$ws.Range("E12").Value = "미래"
$ws.Range("F12").Formula = '="2024-05-23"'
$ws.Range("G12").Formula = '="2024-05-29"'
$ws.Range("H12").Value = 13300
$ws.Range("I12").Value = 6650000
$ws.Range("J12").Value = 2000
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 100

# Row 13: 미래에셋비전스팩6호
$ws.Range("A13").Value = "미래"
$ws.Range("B13").Formula = '="2024-06-13"'
$ws.Range("C13").Value = "미래에셋비전스팩6호"
$ws.Range("D13").Value = "미래"
$ws.Range("E13").Value = "미래"
$ws.Range("F13").Formula = '="2024-06-18"'
$ws.Range("G13").Formula = '="2024-06-24"'
$ws.Range("H13").Value = 12900
$ws.Range("I13").Value = 6450000
$ws.Range("J13").Value = 2000
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 100

# Row 14: 미래에셋비전스팩5호
$ws.Range("A14").Value = "미래"
$ws.Range("B14").Formula = '="2024-06-10"'
$ws.Range("C14").Value = "미래에셋비전스팩5호"
$ws.Range("D14").Value = "미래"
$ws.Range("E14").Value = "미래"
$ws.Range("F14").Formula = '="2024-06-13"'
$ws.Range("G14").Formula = '="2024-06-19"'
$ws.Range("H14").Value = 9500
$ws.Range("I14").Value = 4750000
$ws.Range("J14").Value = 2000
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 100

# Row 15: 하스
$ws.Range("A15").Value = "삼성"
$ws.Range("B15").Formula = '="2024-06-24"'
$ws.Range("C15").Value = "하스"
$ws.Range("D15").Value = "삼성"
$ws.Range("E15").Value = "삼성"
$ws.Range("F15").Formula = '="2024-06-27"'
$ws.Range("G15").Formula = '="2024-07-03"'
$ws.Range("H15").Value = 28960
$ws.Range("I15").Value = 1810000
$ws.Range("J15").Value = 16000
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 100

# Row 16: 그리드위즈
$ws.Range("A16").Value = "삼성"
$ws.Range("B16").Formula = '="2024-06-03"'
$ws.Range("C16").Value = "그리드위즈"
$ws.Range("D16").Value = "삼성"
$ws.Range("E16").Value = "삼성"
$ws.Range("F16").Formula = '="2024-06-07"'
$ws.Range("G16").Formula = '="2024-06-14"'
$ws.Range("H16").Value = 56000
$ws.Range("I16").Value = 1400000
$ws.Range("J16").Value = 40000
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 100

# Row 17: 노브랜드
$ws.Range("A17").Value = "삼성"
$ws.Range("B17").Formula = '="2024-05-13"'
$ws.Range("C17").Value = "노브랜드"
$ws.Range("D17").Value = "삼성"
$ws.Range("E17").Value = "삼성"
$ws.Range("F17").Formula = '="2024-05-17"'
$ws.Range("G17").Formula = '="2024-05-23"'
$ws.Range("H17").Value = 16800
$ws.Range("I17").Value = 1200000
$ws.Range("J17").Value = 14000
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 100

# Row 18: 신한글로벌액티브리츠
$ws.Range("A18").Value = "신한"
$ws.Range("B18").Formula = '="2024-06-13"'
$ws.Range("C18").Value = "신한글로벌액티브리츠"
$ws.Range("D18").Value = "신한, 한국"
$ws.Range("E18").Value = "신한, 한국"
$ws.Range("F18").Formula = '="2024-06-18"'
$ws.Range("G18").Formula = '="2024-07-01"'
$ws.Range("H18").Value = 35000.001
$ws.Range("I18").Value = 23333334
$ws.Range("J18").Value = 3000
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 50

# Row 19: 이노스페이스
$ws.Range("A19").Value = "신한"
$ws.Range("B19").Formula = '="2024-06-20"'
$ws.Range("C19").Value = "이노스페이스"
$ws.Range("D19").Value = "미래"
$ws.Range("E19").Value = "미래, 신한"
$ws.Range("F19").Formula = '="2024-06-25"'
$ws.Range("G19").Formula = '="2024-07-02"'
$ws.Range("H19").Value = 3455.34
$ws.Range("I19").Value = 1330000
$ws.Range("J19").Value = 43300
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 6

# Row 20: 시프트업
$ws.Range("A20").Value = "신한"
$ws.Range("B20").Formula = '="2024-07-02"'
$ws.Range("C20").Value = "시프트업"
$ws.Range("D20").Value = "한국, 제이피모간회사 서울지점, NH"
$ws.Range("E20").Value = "한국, 제이피모간회사 서울지점, NH, 신한"
$ws.Range("F20").Formula = '="2024-07-05"'
$ws.Range("G20").Formula = '="2024-07-11"'
$ws.Range("H20").Value = 17400
$ws.Range("I20").Value = 7250000
$ws.Range("J20").Value = 60000
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 4

# Row 21: 시프트업
$ws.Range("A21").Value = "제이피모간회사"
$ws.Range("B21").Formula = '="2024-07-02"'
$ws.Range("C21").Value = "시프트업"
$ws.Range("D21").Value = "한국, 제이피모간회사 서울지점, NH"
$ws.Range("E21").Value = "한국, 제이피모간회사 서울지점, NH, 신한"
$ws.Range("F21").Formula = '="2024-07-05"'
$ws.Range("G21").Formula = '="2024-07-11"'
$ws.Range("H21").Value = 143550
$ws.Range("I21").Value = 7250000
$ws.Range("J21").Value = 60000
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 33

# Row 22: 한국제15호스팩
$ws.Range("A22").Value = "한국"
$ws.Range("B22").Formula = '="2024-06-17"'
$ws.Range("C22").Value = "한국제15호스팩"
$ws.Range("D22").Value = "한국"
$ws.Range("E22").Value = "한국"
$ws.Range("F22").Formula = '="2024-06-20"'
$ws.Range("G22").Formula = '="2024-06-26"'
$ws.Range("H22").Value = 12500
$ws.Range("I22").Value = 6250000
$ws.Range("J22").Value = 2000
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 100

# Row 23: 한국제14호스팩
$ws.Range("A23").Value = "한국"
$ws.Range("B23").Formula = '="2024-06-10"'
$ws.Range("C23").Value = "한국제14호스팩"
$ws.Range("D23").Value = "한국"
$ws.Range("E23").Value = "한국"
$ws.Range("F23").Formula = '="2024-06-13"'
$ws.Range("G23").Formula = '="2024-06-19"'
$ws.Range("H23").Value = 8000
$ws.Range("I23").Value = 4000000
$ws.Range("J23").Value = 2000
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 100

# Row 24: 씨어스테크놀로지
$ws.Range("A24").Value = "한국"
$ws.Range("B24").Formula = '="2024-06-10"'
$ws.Range("C24").Value = "씨어스테크놀로지"
$ws.Range("D24").Value = "한국"
$ws.Range("E24").Value = "한국"
$ws.Range("F24").Formula = '="2024-06-13"'
$ws.Range("G24").Formula = '="2024-06-19"'
$ws.Range("H24").Value = 22100
$ws.Range("I24").Value = 1300000
$ws.Range("J24").Value = 17000
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 100

# Row 25: 시프트업
$ws.Range("A25").Value = "한국"
$ws.Range("B25").Formula = '="2024-07-02"'
$ws.Range("C25").Value = "시프트업"
$ws.Range("D25").Value = "한국, 제이피모간회사 서울지점, NH"
$ws.Range("E25").Value = "한국, 제이피모간회사 서울지점, NH, 신한"
$ws.Range("F25").Formula = '="2024-07-05"'
$ws.Range("G25").Formula = '="2024-07-11"'
$ws.Range("H25").Value = 143550
$ws.Range("I25").Value = 7250000
$ws.Range("J25").Value = 60000
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 33

# Row 26: 하이젠알앤엠
$ws.Range("A26").Value = "한국"
$ws.Range("B26").Formula = '="2024-06-18"'
$ws.Range("C26").Value = "하이젠알앤엠"
$ws.Range("D26").Value = "한국"
$ws.Range("E26").Value = "한국"
$ws.Range("F26").Formula = '="2024-06-21"'
$ws.Range("G26").Formula = '="2024-06-27"'
$ws.Range("H26").Value = 23800
$ws.Range("I26").Value = 3400000
$ws.Range("J26").Value = 7000
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 100

# Row 27: 에스오에스랩
$ws.Range("A27").Value = "한국"
$ws.Range("B27").Formula = '="2024-06-14"'
$ws.Range("C27").Value = "에스오에스랩"
$ws.Range("D27").Value = "한국"
$ws.Range("E27").Value = "한국, BNK"
$ws.Range("F27").Formula = '="2024-06-19"'
$ws.Range("G27").Formula = '="2024-06-25"'
$ws.Range("H27").Value = 21850
$ws.Range("I27").Value = 2000000
$ws.Range("J27").Value = 11500
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 95

# Row 28: 신한글로벌액티브리츠
$ws.Range("A28").Value = "한국"
$ws.Range("B28").Formula = '="2024-06-13"'
$ws.Range("C28").Value = "신한글로벌액티브리츠"
$ws.Range("D28").Value = "신한, 한국"
$ws.Range("E28").Value = "신한, 한국"
$ws.Range("F28").Formula = '="2024-06-18"'
$ws.Range("G28").Formula = '="2024-07-01"'
$ws.Range("H28").Value = 35000.001
$ws.Range("I28").Value = 23333334
$ws.Range("J28").Value = 3000
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 50

# Row 29: 에이치엠씨제7호스팩
$ws.Range("A29").Value = "현대차"
$ws.Range("B29").Formula = '="2024-06-11"'
$ws.Range("C29").Value = "에이치엠씨제7호스팩"
$ws.Range("D29").Value = "현대차"
$ws.Range("E29").Value = "현대차"
$ws.Range("F29").Formula = '="2024-06-14"'
$ws.Range("G29").Formula = '="2024-06-24"'
$ws.Range("H29").Value = 14000
$ws.Range("I29").Value = 7000000
$ws.Range("J29").Value = 2000
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 100

# The date-like columns (청약일/납입일/상장일) were written as literal-text formulas
# above so Excel does not auto-convert "YYYY-MM-DD" strings into date serials.
# Convert each one back to a plain static text value (matches the source export,
# which stores these as shared-string text, not dates) without touching styles.
$dateCells = @("B2","F2","G2","B3","F3","G3","B4","F4","G4","B5","F5","G5","B6","F6","G6","B7","F7","G7","B8","F8","G8","B9","F9","G9","B10","F10","G10","B11","F11","G11","B12","F12","G12","B13","F13","G13","B14","F14","G14","B15","F15","G15","B16","F16","G16","B17","F17","G17","B18","F18","G18","B19","F19","G19","B20","F20","G20","B21","F21","G21","B22","F22","G22","B23","F23","G23","B24","F24","G24","B25","F25","G25","B26","F26","G26","B27","F27","G27","B28","F28","G28","B29","F29","G29")
foreach ($ref in $dateCells) {
    $ws.Range($ref).Copy() | Out-Null
    $ws.Range($ref).PasteSpecial(-4163) | Out-Null
}
$excel.CutCopyMode = 0
